# Natmi following Dr Hou advice
# Rewrites the Cd80-Cd28 LR-pair sheet: sending clusters now include "sCs"
# (5 sending clusters x 3 target clusters = rows 2-16) with refreshed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd80"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6166306666666667
$ws.Range("H2").Value = 1.849892
$ws.Range("I2").Value = 0.01585912217154475
$ws.Range("J2").Value = 0.01594598818140205
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.473964666666667
$ws.Range("N2").Value = 10.421894
$ws.Range("O2").Value = 0.5449432418219904
$ws.Range("P2").Value = 0.5462679851731443
$ws.Range("Q2").Value = 2.142153148383111
$ws.Range("R2").Value = 19.279378335448
$ws.Range("S2").Value = 0.008642321448612598
$ws.Range("T2").Value = 0.008710782835449269

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd80"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6166306666666667
$ws.Range("H3").Value = 1.849892
$ws.Range("I3").Value = 0.01585912217154475
$ws.Range("J3").Value = 0.01594598818140205
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.854567333333333
$ws.Range("N3").Value = 8.563701999999999
$ws.Range("O3").Value = 0.4477815193550676
$ws.Range("P3").Value = 0.4488700649961731
$ws.Range("Q3").Value = 1.760213757798222
$ws.Range("R3").Value = 15.841923820184
$ws.Range("S3").Value = 0.007101421821611946
$ws.Range("T3").Value = 0.007157676751414145

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cd80"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6166306666666667
$ws.Range("H4").Value = 1.849892
$ws.Range("I4").Value = 0.01585912217154475
$ws.Range("J4").Value = 0.01594598818140205
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.046379
$ws.Range("N4").Value = 0.09275799999999999
$ws.Range("O4").Value = 0.007275238822941998
$ws.Range("P4").Value = 0.004861949830682458
$ws.Range("Q4").Value = 0.02859871368933333
$ws.Range("R4").Value = 0.171592282136
$ws.Range("S4").Value = 0.0001153789013202025
$ws.Range("T4").Value = 0.00007752859453863216

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd80"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "M1"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.510701
$ws.Range("H5").Value = 7.532103
$ws.Range("I5").Value = 0.06457271110186903
$ws.Range("J5").Value = 0.06492639863251634
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.473964666666667
$ws.Range("N5").Value = 10.421894
$ws.Range("O5").Value = 0.5449432418219904
$ws.Range("P5").Value = 0.5462679851731443
$ws.Range("Q5").Value = 8.722086562564668
$ws.Range("R5").Value = 78.498779063082
$ws.Range("S5").Value = 0.03518846252108734
$ws.Range("T5").Value = 0.0354672129655331

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cd80"
$ws.Range("C6").Value = "Cd28"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.510701
$ws.Range("H6").Value = 7.532103
$ws.Range("I6").Value = 0.06457271110186903
$ws.Range("J6").Value = 0.06492639863251634
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.854567333333333
$ws.Range("N6").Value = 8.563701999999999
$ws.Range("O6").Value = 0.4477815193550676
$ws.Range("P6").Value = 0.4488700649961731
$ws.Range("Q6").Value = 7.166965058367333
$ws.Range("R6").Value = 64.50268552530599
$ws.Range("S6").Value = 0.02891446668607076
$ws.Range("T6").Value = 0.02914351677414505

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cd80"
$ws.Range("C7").Value = "Cd28"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.510701
$ws.Range("H7").Value = 7.532103
$ws.Range("I7").Value = 0.06457271110186903
$ws.Range("J7").Value = 0.06492639863251634
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.046379
$ws.Range("N7").Value = 0.09275799999999999
$ws.Range("O7").Value = 0.007275238822941998
$ws.Range("P7").Value = 0.004861949830682458
$ws.Range("Q7").Value = 0.116443801679
$ws.Range("R7").Value = 0.698662810074
$ws.Range("S7").Value = 0.0004697818947109353
$ws.Range("T7").Value = 0.0003156688928381846

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Cd80"
$ws.Range("C8").Value = "Cd28"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.23029433333333
$ws.Range("H8").Value = 36.690883
$ws.Range("I8").Value = 0.314550901392543
$ws.Range("J8").Value = 0.3162738077051013
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.473964666666667
$ws.Range("N8").Value = 10.421894
$ws.Range("O8").Value = 0.5449432418219904
$ws.Range("P8").Value = 0.5462679851731443
$ws.Range("Q8").Value = 42.48761037693356
$ws.Range("R8").Value = 382.388493392402
$ws.Range("S8").Value = 0.1714123879228816
$ws.Range("T8").Value = 0.1727702556981042

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Cd80"
$ws.Range("C9").Value = "Cd28"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.23029433333333
$ws.Range("H9").Value = 36.690883
$ws.Range("I9").Value = 0.314550901392543
$ws.Range("J9").Value = 0.3162738077051013
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.854567333333333
$ws.Range("N9").Value = 8.563701999999999
$ws.Range("O9").Value = 0.4477815193550676
$ws.Range("P9").Value = 0.4488700649961731
$ws.Range("Q9").Value = 34.91219868098511
$ws.Range("R9").Value = 314.209788128866
$ws.Range("S9").Value = 0.140850080540059
$ws.Range("T9").Value = 0.141965844621176

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Cd80"
$ws.Range("C10").Value = "Cd28"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.23029433333333
$ws.Range("H10").Value = 36.690883
$ws.Range("I10").Value = 0.314550901392543
$ws.Range("J10").Value = 0.3162738077051013
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.046379
$ws.Range("N10").Value = 0.09275799999999999
$ws.Range("O10").Value = 0.007275238822941998
$ws.Range("P10").Value = 0.004861949830682458
$ws.Range("Q10").Value = 0.5672288208856666
$ws.Range("R10").Value = 3.403372925314
$ws.Range("S10").Value = 0.002288432929602429
$ws.Range("T10").Value = 0.001537707385821113

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Cd80"
$ws.Range("C11").Value = "Cd28"
$ws.Range("D11").Value = "M1"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 22.88871233333333
$ws.Range("H11").Value = 68.66613700000001
$ws.Range("I11").Value = 0.5886747203247699
$ws.Range("J11").Value = 0.591899099549884
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.473964666666667
$ws.Range("N11").Value = 10.421894
$ws.Range("O11").Value = 0.5449432418219904
$ws.Range("P11").Value = 0.5462679851731443
$ws.Range("Q11").Value = 79.51457791149757
$ws.Range("R11").Value = 715.6312012034781
$ws.Range("S11").Value = 0.3207943104724336
$ws.Range("T11").Value = 0.3233355285369135

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Cd80"
$ws.Range("C12").Value = "Cd28"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 22.88871233333333
$ws.Range("H12").Value = 68.66613700000001
$ws.Range("I12").Value = 0.5886747203247699
$ws.Range("J12").Value = 0.591899099549884
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.854567333333333
$ws.Range("N12").Value = 8.563701999999999
$ws.Range("O12").Value = 0.4477815193550676
$ws.Range("P12").Value = 0.4488700649961731
$ws.Range("Q12").Value = 65.33737052879711
$ws.Range("R12").Value = 588.036334759174
$ws.Range("S12").Value = 0.2635976606729449
$ws.Range("T12").Value = 0.2656857872861327

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Cd80"
$ws.Range("C13").Value = "Cd28"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 22.88871233333333
$ws.Range("H13").Value = 68.66613700000001
$ws.Range("I13").Value = 0.5886747203247699
$ws.Range("J13").Value = 0.591899099549884
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.046379
$ws.Range("N13").Value = 0.09275799999999999
$ws.Range("O13").Value = 0.007275238822941998
$ws.Range("P13").Value = 0.004861949830682458
$ws.Range("Q13").Value = 1.061555589307667
$ws.Range("R13").Value = 6.369333535846
$ws.Range("S13").Value = 0.004282749179391289
$ws.Range("T13").Value = 0.002877783726837657

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Cd80"
$ws.Range("C14").Value = "Cd28"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 0.635427
$ws.Range("H14").Value = 1.270854
$ws.Range("I14").Value = 0.01634254500927324
$ws.Range("J14").Value = 0.01095470593109626
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.473964666666667
$ws.Range("N14").Value = 10.421894
$ws.Range("O14").Value = 0.5449432418219904
$ws.Range("P14").Value = 0.5462679851731443
$ws.Range("Q14").Value = 2.207450946246
$ws.Range("R14").Value = 13.244705677476
$ws.Range("S14").Value = 0.00890575945697515
$ws.Range("T14").Value = 0.005984205137144247

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Cd80"
$ws.Range("C15").Value = "Cd28"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.5
$ws.Range("G15").Value = 0.635427
$ws.Range("H15").Value = 1.270854
$ws.Range("I15").Value = 0.01634254500927324
$ws.Range("J15").Value = 0.01095470593109626
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.854567333333333
$ws.Range("N15").Value = 8.563701999999999
$ws.Range("O15").Value = 0.4477815193550676
$ws.Range("P15").Value = 0.4488700649961731
$ws.Range("Q15").Value = 1.813869156918
$ws.Range("R15").Value = 10.883214941508
$ws.Range("S15").Value = 0.00731788963438095
$ws.Range("T15").Value = 0.00491723956330514

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Cd80"
$ws.Range("C16").Value = "Cd28"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.5
$ws.Range("G16").Value = 0.635427
$ws.Range("H16").Value = 1.270854
$ws.Range("I16").Value = 0.01634254500927324
$ws.Range("J16").Value = 0.01095470593109626
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.5
$ws.Range("M16").Value = 0.046379
$ws.Range("N16").Value = 0.09275799999999999
$ws.Range("O16").Value = 0.007275238822941998
$ws.Range("P16").Value = 0.004861949830682458
$ws.Range("Q16").Value = 0.029470468833
$ws.Range("R16").Value = 0.117881875332
$ws.Range("S16").Value = 0.0001188959179171417
$ws.Range("T16").Value = 0.00005326123064686957

